$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()
[void]$ws.Rows.Item(16).EntireRow.Select()
